$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (the GUSTAVO/004565108 row) and populate it
# with the RODRIGO/004392159/900.21 record. Force the account-number cell
# to Text format first so the leading zeros in "004392159" survive (Excel
# would otherwise interpret the literal as a number).
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004392159"
$ws.Cells.Item(4, 2).Value = "RODRIGO"
$ws.Cells.Item(4, 3).Value = 900.21

# The GUSTAVO row (now shifted down to row 5) gets its balance corrected.
$ws.Cells.Item(5, 3).Value = 432.86

# Remove the now-redundant rows that used to follow GUSTAVO: the duplicate
# RODRIGO row plus the LAURA, BLUEMETRIX, LEIVANIO and IRON rows (rows 6-10
# after the insert above).
$ws.Range("A6:A10").EntireRow.Delete()
